$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at position 5 (pushes the old blank formatted row 5 down to row 7,
# and copies the row-4 formatting into the two new rows 5 and 6)
$ws.Range("A5:A6").EntireRow.Insert()

# Common values shared across all data rows (only ID_Building/Hop/Htr_w vary by row)
$type = "SFH"
$constructionStart = 1949
$constructionEnd = 1957
$personNum = 3
$af = 174.14
$hve = 49.261000000000003
$cmFactor = 213505.516
$amFactor = 3
$internalGains = 3
$effWestEast = 10.037000000000001
$effSouth = 4.5330000000000004
$effNorth = 1.619
$gridPowerMax = 21000

# Per-row values: ID_Building (A), Hop (G), Htr_w (H)
$rowData = @(
    @{ Row = 2; Id = 1; Hop = 303.87599999999998; HtrW = 117.355 },
    @{ Row = 3; Id = 2; Hop = 286.52800000000002; HtrW = 115.621 },
    @{ Row = 4; Id = 3; Hop = 274.81700000000001; HtrW = 114.45 },
    @{ Row = 5; Id = 4; Hop = 273.15499999999997; HtrW = 114.283 },
    @{ Row = 6; Id = 5; Hop = 303.87599999999998; HtrW = 136.28700000000001 }
)

foreach ($rd in $rowData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $rd.Id
    $ws.Cells.Item($r, 2).Value = $type
    $ws.Cells.Item($r, 3).Value = $constructionStart
    $ws.Cells.Item($r, 4).Value = $constructionEnd
    $ws.Cells.Item($r, 5).Value = $personNum
    $ws.Cells.Item($r, 6).Value = $af
    $ws.Cells.Item($r, 7).Value = $rd.Hop
    $ws.Cells.Item($r, 8).Value = $rd.HtrW
    $ws.Cells.Item($r, 9).Value = $hve
    $ws.Cells.Item($r, 10).Value = $cmFactor
    $ws.Cells.Item($r, 11).Value = $amFactor
    $ws.Cells.Item($r, 12).Value = $internalGains
    $ws.Cells.Item($r, 13).Value = $effWestEast
    $ws.Cells.Item($r, 14).Value = $effSouth
    $ws.Cells.Item($r, 15).Value = $effNorth
    $ws.Cells.Item($r, 16).Value = $gridPowerMax
}

# Make sure the two freshly-inserted rows carry the same explicit cell style as the rest
# of the data block (re-asserting an automatic font color is a no-op visually but forces
# the style index to be written out on every cell).
$ws.Range("A5:P6").Font.ColorIndex = -4105

# Nudge the worksheet's used range so it still reaches row 7 (the formatted-but-empty
# row that used to be row 5) the same way it covered row 5 before this edit.
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 1).ClearContents()

# Refresh the autofilter so its range covers the newly added rows. Toggling it off first
# is required because calling AutoFilter() while a filter is already active would simply
# remove it instead of re-applying it with the new range.
$ws.AutoFilterMode = $false
$ws.Range("A1:P6").AutoFilter()

# The _FilterDatabase defined name is not refreshed automatically by the interop layer,
# so update it by hand to track the new autofilter extent.
foreach ($n in $ws.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = '=OperationScenario_Component_Bui!$A$1:$P$6'
    }
}

# Match the post-edit selection recorded in the workbook
$ws.Range("C10").Select()
